$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Range("E1").Value = "BirthDate"

# Date+time value (serial date 34100.11037037037 == 1993-05-11 02:38:56).
# Assign the raw serial number (not a DateTime) so Excel doesn't silently
# stamp on its own default datetime display format - we set the exact
# builtin number format ourselves below.
$birthDate = 34100.11037037037

$ws.Range("E2").Value = $birthDate
$ws.Range("E2").NumberFormat = "m/d/yy h:mm"

# Row 3 column E is intentionally left empty (reading the last empty cell of a row)

$ws.Range("E4").Value = $birthDate
$ws.Range("E4").NumberFormat = "m/d/yy h:mm"

# Resize column E to fit the new content (matches the author's best-fit column width)
$ws.Columns.Item(5).ColumnWidth = 13.3

# Move the active selection, matching the author's final cursor position
$ws.Range("F8").Select()
